# HungryDragonContent_Gameplay.xlsx - "added sku for AnniversaryCandle"
# Adds a new entity definition row (AnniversaryCandle) right after the
# existing "Anniversary1_Cake_Piece" row in the entityDefinitions table,
# tweaks that existing row's reward numbers, re-applies the [rewardScore]
# AutoFilter (which hides a handful of out-of-range rows), and keeps the
# tables/defined names that live below the insertion point in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entities")

# ---------------------------------------------------------------------
# 1) Tweak the existing "Anniversary1_Cake_Piece" row (row 157) rewards
# ---------------------------------------------------------------------
$ws.Range("D157").Value = 40
$ws.Range("E157").Value = 12
$ws.Range("G157").Value = 30

# ---------------------------------------------------------------------
# 2) Insert a new row right below it for the new "AnniversaryCandle" sku
#    and copy row 157's formatting down so the new row looks the same.
# ---------------------------------------------------------------------
$ws.Range("A157:AF157").Copy()
$ws.Rows.Item(158).Insert(-4121)
$excel.CutCopyMode = $false

$ws.Range("A158").Value = "<Definition>"
$ws.Range("B158").Value = "AnniversaryCandle"
$ws.Range("C158").Value = "collectible"
$ws.Range("D158").Value = 20
$ws.Range("E158").Value = 0
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0
$ws.Range("I158").Value = 25
$ws.Range("J158").Value = 0.1
$ws.Range("K158").Value = 0
$ws.Range("L158").Value = 0
$ws.Range("M158").Value = $false
$ws.Range("N158").Value = 5
$ws.Range("O158").Value = 5
$ws.Range("P158").Value = 0
$ws.Range("Q158").Value = 0
$ws.Range("R158").Value = $true
$ws.Range("S158").Value = $false
$ws.Range("T158").Value = $false
$ws.Range("U158").Value = 1
$ws.Range("V158").Value = ""
$ws.Range("W158").Value = 0
$ws.Range("X158").Value = 0
$ws.Range("Y158").Value = 0
$ws.Range("Z158").Value = 0
$ws.Range("AA158").Value = 0
$ws.Range("AB158").Value = 0
$ws.Range("AC158").Value = "TID_QUIP_DRG_KILL_ENT_BIRTHDAY_CAKE_02"
$ws.Range("AD158").Value = "TID_QUIP_DRG_BURN_ENT_09"
$ws.Range("AE158").Value = ""
$ws.Range("AF158").Value = ""

# ---------------------------------------------------------------------
# 3) Grow the entityDefinitions table to include the new row, and shift
#    the decorationDefinitions ("Table4") table down by one row so it
#    still lines up with its (now shifted) header/data rows.
# ---------------------------------------------------------------------
$entityTable = $ws.ListObjects.Item("entityDefinitions")
$entityTable.Resize($ws.Range("A23:AF158"))

$decoTable = $ws.ListObjects.Item("Table4")
$decoTable.Resize($ws.Range("A163:O175"))

# ---------------------------------------------------------------------
# 4) Re-apply the [rewardScore] column AutoFilter on entityDefinitions.
#    This both writes the filter criteria back to the table and hides
#    the rows whose [rewardScore] isn't one of the kept values.
# ---------------------------------------------------------------------
$rewardScoreValues = @(100,110,120,180,20,220,270,30,330,360,40,400,50,500,540,60,65,80,810,85,90)
$entityTable.Range.AutoFilter(4, $rewardScoreValues, 7)

# ---------------------------------------------------------------------
# 5) Fix up the hidden _xlnm._FilterDatabase defined name so it still
#    points at the (now shifted down by one row) Table4 header+row.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -eq "entities!_FilterDatabase") {
    $n.RefersTo = "=entities!`$A`$163:`$O`$164"
  }
}

# ---------------------------------------------------------------------
# 6) Leave the view looking at the newly added row, like the author did.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("N158").Select()
